$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column G ("K" = strikeouts) values per regenerated save_data
$ws.Range("G2").Value = 5
$ws.Range("G3").Value = 3
$ws.Range("G4").Value = 0
$ws.Range("G5").Value = 0
